$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the whole duty-name column (B2:B32), then re-populate only the
# rows that keep an assigned name, per the new shortened roster.
$ws.Range("B2:B32").ClearContents()

$ws.Range("B2").Value  = "兒島大志郎"
$ws.Range("B6").Value  = "高野怜央"
$ws.Range("B10").Value = "志塚惇希"
$ws.Range("B14").Value = "白岩詩佑介"
$ws.Range("B18").Value = "Nicholas Tristan Aryasatyo"
$ws.Range("B22").Value = "川田涼介"
$ws.Range("B26").Value = "三神佳誠"
$ws.Range("B30").Value = "兒島大志郎"

# B5/B28 previously carried the one-off "Roboto" font used only for the
# now-removed "Ethan Virtudazo" entries; restore them to the standard
# Arial styling shared by the rest of the column now that it's blank.
$ws.Range("B5").Font.Name = "Arial"
$ws.Range("B5").Font.Size = 10
$ws.Range("B28").Font.Name = "Arial"
$ws.Range("B28").Font.Size = 10

# Move the active selection to D31, matching the saved view state.
$ws.Range("D31").Select()
